$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells being updated so numeric-looking strings
# (e.g. "1.001", "25.078.93") are preserved as text, matching the source data.
$priceCells = @("D2","D3","D4","D5","D7","D8","D9","D10","D11","D12","D13","D15","D16","D17","D18","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D41","D42","D43","D44","D45","D46","D47","D48","D49","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '25.078.93'
$ws.Range("E2").Value = '  -2.97%  '

$ws.Range("D3").Value = '1.653.05'
$ws.Range("E3").Value = '  -4.84%  '

$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = '236.69'
$ws.Range("E5").Value = '  -1.70%  '

$ws.Range("E6").Value = '  +0.03%  '

$ws.Range("D7").Value = '0.4783'
$ws.Range("E7").Value = '  -7.72%  '

$ws.Range("D8").Value = '0.2621'
$ws.Range("E8").Value = '  -4.36%  '

$ws.Range("D9").Value = '0.05977'
$ws.Range("E9").Value = '  -2.81%  '

$ws.Range("D10").Value = '0.07110'
$ws.Range("E10").Value = '  -0.79%  '

$ws.Range("D11").Value = '1.658.67'
$ws.Range("E11").Value = '  -4.51%  '

$ws.Range("D12").Value = '14.42'
$ws.Range("E12").Value = '  -3.36%  '

$ws.Range("D13").Value = '0.6197'
$ws.Range("E13").Value = '  -3.31%  '

$ws.Range("E14").Value = '  -0.12%  '

$ws.Range("D15").Value = '73.04'
$ws.Range("E15").Value = '  -5.39%  '

$ws.Range("D16").Value = '1.000'
$ws.Range("E16").Value = '  -0.01%  '

$ws.Range("D17").Value = '1.001'
$ws.Range("E17").Value = '  +0.05%  '

$ws.Range("D18").Value = '25.074.42'
$ws.Range("E18").Value = '  -3.09%  '

$ws.Range("E19").Value = '  -2.95%  '

$ws.Range("D20").Value = '0.000006550'
$ws.Range("E20").Value = '  -3.22%  '

$ws.Range("D21").Value = '4.445'
$ws.Range("E21").Value = '  +4.39%  '

$ws.Range("D22").Value = '1.865.20'
$ws.Range("E22").Value = '  -5.01%  '

$ws.Range("D23").Value = '8.468'
$ws.Range("E23").Value = '  -1.83%  '

$ws.Range("D24").Value = '5.267'
$ws.Range("E24").Value = '  +0.26%  '

$ws.Range("D25").Value = '133.23'
$ws.Range("E25").Value = '  -3.50%  '

$ws.Range("D26").Value = '14.72'
$ws.Range("E26").Value = '  -3.15%  '

$ws.Range("D27").Value = '1.398'
$ws.Range("E27").Value = '  -8.02%  '

$ws.Range("D28").Value = '1.696'

$ws.Range("D29").Value = '101.59'
$ws.Range("E29").Value = '  -3.17%  '

$ws.Range("D30").Value = '3.803'
$ws.Range("E30").Value = '  -3.52%  '

$ws.Range("D31").Value = '0.07904'
$ws.Range("E31").Value = '  -4.27%  '

$ws.Range("D32").Value = '3.522'
$ws.Range("E32").Value = '  -3.76%  '

$ws.Range("D33").Value = '0.04594'
$ws.Range("E33").Value = '  -0.75%  '

$ws.Range("D34").Value = '2.604'
$ws.Range("E34").Value = '  -1.52%  '

$ws.Range("D35").Value = '0.9408'
$ws.Range("E35").Value = '  -4.59%  '

$ws.Range("D36").Value = '0.5841'
$ws.Range("E36").Value = '  -5.08%  '

$ws.Range("D37").Value = '2.622'
$ws.Range("E37").Value = '  -2.63%  '

$ws.Range("D38").Value = '0.01536'
$ws.Range("E38").Value = '  -3.87%  '

$ws.Range("D39").Value = '0.8418'
$ws.Range("E39").Value = '  +13.05%  '

$ws.Range("E40").Value = '  +0.08%  '

$ws.Range("D41").Value = '1.830'
$ws.Range("E41").Value = '  -4.51%  '

$ws.Range("D42").Value = '98.50'
$ws.Range("E42").Value = '  -1.06%  '

$ws.Range("D43").Value = '0.3698'
$ws.Range("E43").Value = '  -3.62%  '

$ws.Range("D44").Value = '4.833'
$ws.Range("E44").Value = '  -3.21%  '

$ws.Range("D45").Value = '0.1117'
$ws.Range("E45").Value = '  -0.56%  '

$ws.Range("D46").Value = '6.050'
$ws.Range("E46").Value = '  -2.95%  '

$ws.Range("D47").Value = '0.05150'
$ws.Range("E47").Value = '  -1.74%  '

$ws.Range("D48").Value = '52.09'
$ws.Range("E48").Value = '  -4.82%  '

$ws.Range("D49").Value = '29.32'
$ws.Range("E49").Value = '  -3.83%  '

$ws.Range("E50").Value = '  -0.07%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '7.306'
$ws.Range("E51").Value = '  -3.81%  '
